# Update gh-pages to output generated at 456a3b4
# This script updates the "报名人数/F" (and one "G") columns across the
# four worksheets of the workbook to reflect newly generated stats.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 68
$ws.Range("F6").Value = 897
$ws.Range("F7").Value = 487
$ws.Range("F8").Value = 4813
$ws.Range("F9").Value = 4813
$ws.Range("F12").Value = 174
$ws.Range("F16").Value = 7769
$ws.Range("F17").Value = 258
$ws.Range("F18").Value = 131
$ws.Range("F21").Value = 1434
$ws.Range("F22").Value = 1434
$ws.Range("F30").Value = 6233
$ws.Range("F31").Value = 154
$ws.Range("F32").Value = 44
$ws.Range("F36").Value = 6585
$ws.Range("F37").Value = 29
$ws.Range("F41").Value = 24
$ws.Range("F47").Value = 44
$ws.Range("F48").Value = 462
$ws.Range("F49").Value = 2172
$ws.Range("G49").Value = 80
$ws.Range("F50").Value = 53

# ---- Sheet "演出" (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 135
$ws.Range("F7").Value = 40

# ---- Sheet "本地生活" (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1460

# ---- Sheet "全部类型" (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1460
$ws.Range("F5").Value = 68
$ws.Range("F10").Value = 487
$ws.Range("F11").Value = 4813
$ws.Range("F12").Value = 4813
$ws.Range("F15").Value = 174
$ws.Range("F17").Value = 7770
$ws.Range("F18").Value = 258
$ws.Range("F19").Value = 131
$ws.Range("F21").Value = 1434
$ws.Range("F22").Value = 135
$ws.Range("F25").Value = 40
$ws.Range("F32").Value = 6233
$ws.Range("F33").Value = 154
$ws.Range("F35").Value = 44
$ws.Range("F37").Value = 6585
$ws.Range("F38").Value = 29
$ws.Range("F42").Value = 24
$ws.Range("F47").Value = 44
$ws.Range("F48").Value = 462
$ws.Range("F50").Value = 53
